# Add new time-log entries to the table on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows to append below the existing data (rows 2-7 already populated).
# Columns: A=Date, B=Duration (hrs), C=Task, D=Comments
$newRows = @(
    @{ Row = 8;  Date = "9/22/2023"; Duration = 0.10416666666666667; Task = "Started working on the app"; Comments = "Started working on the main page with the balance. So far getting the hang of Xcode. Got the top of the page icons complete as well as the logo in the middle."; Height = 48 },
    @{ Row = 9;  Date = "9/24/2023"; Duration = 3; Task = 'Learning how to use "buttons"'; Comments = "Making buttons to bring up other screens in Xcode was more challanging than I thought so I spent quite some time reading and watching videos on how they work."; Height = 48 },
    @{ Row = 10; Date = "9/25/2023"; Duration = 1; Task = "Debugging"; Comments = "Main page kept showing a black recangle on the top of the screen and I couldn't figure out where the Zstack should go in the code."; Height = 36 },
    @{ Row = 11; Date = "9/27/2023"; Duration = 0.0625; Task = "Greeting message"; Comments = "The greeting message on the main screen has to say good morning, afternoon, and evening based on time of day. I was able to figoure out how to make it behave like that."; Height = 48 },
    @{ Row = 12; Date = "9/29/2023"; Duration = 2; Task = "Navigation "; Comments = "I was able to make the buttons be able to take the user to the corresponding page. Finding the fucntion to do that wasn't the issue but it kept resizing everything and I spent a while figuring that out. Turns out I had an extra " + [char]34 + "Spacer()" + [char]34 + " written and an extra " + [char]34 + "}" + [char]34; Height = 72 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Clone the formatting of the last existing data row (row 7) onto the
    # new row so dates/tasks/comments pick up the same styles (s2/s1/s6).
    $ws.Range("A7:D7").Copy()
    $ws.Range("A" + $row + ":D" + $row).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 3).Value = $r.Task
    $ws.Cells.Item($row, 4).Value = $r.Comments

    if ($r.Duration -lt 1) {
        # Sub-hour durations keep the time-of-day format already cloned
        # from row 7 (style index 9 - h:mm, centered).
        $ws.Cells.Item($row, 2).Value = $r.Duration
    } else {
        # Whole-hour durations use the plain centered/wrapped General
        # format (style index 3). Start from an untouched cell's format
        # (General, no explicit number format) so the engine reuses the
        # existing style instead of minting a new numFmt record.
        $ws.Range("Z1").Copy()
        $ws.Cells.Item($row, 2).PasteSpecial(-4122)
        $ws.Cells.Item($row, 2).Value = $r.Duration
        $ws.Cells.Item($row, 2).HorizontalAlignment = -4108
        $ws.Cells.Item($row, 2).VerticalAlignment = -4108
        $ws.Cells.Item($row, 2).WrapText = $true
    }

    $ws.Rows.Item($row).RowHeight = $r.Height
}

# Grow the table to include the newly added rows.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:D12"))

# Reflect the selection left behind after data entry (last empty row).
[void]$ws.Range("A13").Select()
